$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7 ("Manuela" -> "Fernanda", senha 123 -> 1101, extrato 8278 -> 10278)
$ws.Range("B7").Value = "Fernanda"
$ws.Range("D7").Value = 1101
$ws.Range("E7").Value = 10278

# Row 2 extrato: 510 -> 10
$ws.Range("E2").Value = 10

# Row 3 extrato: 1092.91 -> 1192.91
$ws.Range("E3").Value = 1192.91

# Update the active selection to E2
$ws.Range("E2").Select()
